$wb = $excel.ActiveWorkbook

# --- Sheet: CypherOutput_Message (copy of Message sheet content) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCOM = $wb.Worksheets.Add($null, $lastSheet)
$wsCOM.Name = "CypherOutput_Message"
$wsCOM.Range("A1").Value = 'Neo4j_URL:'
$wsCOM.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsCOM.Range("A3").Value = 'User_name:'
$wsCOM.Range("A4").Value = 'neo4j'
$wsCOM.Range("A5").Value = 'PWD:'
$wsCOM.Range("A6").Value = 'icdcDBneo4j0'
$wsCOM.Range("A7").Value = 'Cypher:'
$wsCOM.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (single lung lobe)''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$wsCOM.Range("A9").Value = 'Output:'
$wsCOM.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC13_Canine_Filter_Diagnosis-MaligSingle_Neo4jData.xlsx'

# --- Sheet: StatOutput (summary counts) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStat = $wb.Worksheets.Add($null, $lastSheet)
$wsStat.Name = "StatOutput"
$wsStat.Range("A1").Value = "number_of_files"
$wsStat.Range("B1").Value = "number_of_sample"
$wsStat.Range("C1").Value = "number_of_cases"
$wsStat.Range("D1").Value = "number_of_study"
$wsStat.Range("A2").Value = "'9"
$wsStat.Range("B2").Value = "'19"
$wsStat.Range("C2").Value = "'9"
$wsStat.Range("D2").Value = "'1"

# --- Sheet: StatOutput_Message (copy of Message sheet content, twice, second with StatOutput cypher) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSOM = $wb.Worksheets.Add($null, $lastSheet)
$wsSOM.Name = "StatOutput_Message"
$wsSOM.Range("A1").Value = 'Neo4j_URL:'
$wsSOM.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsSOM.Range("A3").Value = 'User_name:'
$wsSOM.Range("A4").Value = 'neo4j'
$wsSOM.Range("A5").Value = 'PWD:'
$wsSOM.Range("A6").Value = 'icdcDBneo4j0'
$wsSOM.Range("A7").Value = 'Cypher:'
$wsSOM.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (single lung lobe)''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$wsSOM.Range("A9").Value = 'Output:'
$wsSOM.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC13_Canine_Filter_Diagnosis-MaligSingle_Neo4jData.xlsx'
$wsSOM.Range("A11").Value = 'Neo4j_URL:'
$wsSOM.Range("A12").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsSOM.Range("A13").Value = 'User_name:'
$wsSOM.Range("A14").Value = 'neo4j'
$wsSOM.Range("A15").Value = 'PWD:'
$wsSOM.Range("A16").Value = 'icdcDBneo4j0'
$wsSOM.Range("A17").Value = 'Cypher:'
$wsSOM.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (single lung lobe)'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$wsSOM.Range("A19").Value = 'Output:'
$wsSOM.Range("A20").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC13_Canine_Filter_Diagnosis-MaligSingle_Neo4jData.xlsx'

# Restore original active sheet/selection so we don't leave tabSelected on a new sheet
$wb.Worksheets.Item("CypherOutput").Select()
